# Auto-generated edit script: appends new daily snapshot rows
# to each of the three portfolio sheets, and refreshes the
# worksheet dimension (handled automatically by the engine).

$wb = $excel.ActiveWorkbook

# --- Worksheet 1: rows 116-126 ---
$ws = $wb.Worksheets.Item(1)
$rows = @(
    @("大智 (稳健智远)", "000333", "美的集团", 3.04, 42.89719527444501, 71.72, 3076.586845083196, 101331.2684224588, "202506231600"),
    @("大智 (稳健智远)", "510050", "上证50ETF", 5.05, 1852.638869852698, 2.76, 5113.283280793446, 101331.2684224588, "202506231600"),
    @("大智 (稳健智远)", "510300", "沪深300ETF", 4.9, 1276.881426590205, 3.89, 4967.068749435897, 101331.2684224588, "202506231600"),
    @("大智 (稳健智远)", "518880", "黄金ETF", 4.96, 673.9096418114974, 7.46, 5027.365927913771, 101331.2684224588, "202506231600"),
    @("大智 (稳健智远)", "600085", "同仁堂", 1.99, 56.06334956913254, 36, 2018.280584488771, 101331.2684224588, "202506231600"),
    @("大智 (稳健智远)", "600900", "长江电力", 20.24, 669.0422707938175, 30.65, 20506.1455998305, 101331.2684224588, "202506231600"),
    @("大智 (稳健智远)", "600989", "宝丰能源", 5.05, 308.5861230826723, 16.58, 5116.357920710707, 101331.2684224588, "202506231600"),
    @("大智 (稳健智远)", "601899", "XD紫金矿", 9.890000000000001, 541.4194359293219, 18.51, 10021.67375905175, 101331.2684224588, "202506231600"),
    @("大智 (稳健智远)", "HK02899", "紫金矿业", 9.91, 521.4694874201556, 19.26, 10043.5023277122, 101331.2684224588, "202506231600"),
    @("大智 (稳健智远)", "HK06881", "中国银河", 4.81, 600.7968033130801, 8.109999999999999, 4872.46207486908, 101331.2684224588, "202506231600"),
    @("大智 (稳健智远)", "100000", "现金", 30.17, 30568.54135256952, 1, 30568.54135256952, 101331.2684224588, "202506231600")
)
$startRow = 116
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = "'" + $row[8]
}

# --- Worksheet 2: rows 79-85 ---
$ws = $wb.Worksheets.Item(2)
$rows = @(
    @("大成 (锐进先锋)", "000725", "京东方A", 5.31, 1264.047242468624, 3.96, 5005.627080175751, 94292.632964472, "202506231600"),
    @("大成 (锐进先锋)", "159781", "科创创业ETF", 5.21, 9277.629760760277, 0.53, 4917.143773202947, 94292.632964472, "202506231600"),
    @("大成 (锐进先锋)", "513100", "纳指ETF", 5.18, 3131.938709046463, 1.56, 4885.824386112482, 94292.632964472, "202506231600"),
    @("大成 (锐进先锋)", "513290", "纳指生物科技ETF", 1.02, 870.2909333102562, 1.1, 957.3200266412819, 94292.632964472, "202506231600"),
    @("大成 (锐进先锋)", "603119", "浙江荣泰", 44.27, 1069.978093782073, 39.01, 41739.84543843867, 94292.632964472, "202506231600"),
    @("大成 (锐进先锋)", "688290", "景业智能", 7.72, 147.121987895035, 49.51, 7284.009620683183, 94292.632964472, "202506231600"),
    @("大成 (锐进先锋)", "100000", "现金", 31.29, 29502.86263921768, 1, 29502.86263921768, 94292.632964472, "202506231600")
)
$startRow = 79
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = "'" + $row[8]
}

# --- Worksheet 3: rows 166-180 ---
$ws = $wb.Worksheets.Item(3)
$rows = @(
    @("范式进化投资组合", "000333", "美的集团", 1.01, 14.02515882310653, 71.72, 1005.8843907932, 99639.24647381395, "202506231600"),
    @("范式进化投资组合", "000725", "京东方A", 5.11, 1284.618017653691, 3.96, 5087.087349908616, 99639.24647381395, "202506231600"),
    @("范式进化投资组合", "159781", "科创创业ETF", 5.02, 9428.611488061997, 0.53, 4997.164088672858, 99639.24647381395, "202506231600"),
    @("范式进化投资组合", "510050", "上证50ETF", 5.03, 1817.150577699221, 2.76, 5015.33559444985, 99639.24647381395, "202506231600"),
    @("范式进化投资组合", "510300", "沪深300ETF", 4.89, 1252.422077361618, 3.89, 4871.921880936694, 99639.24647381395, "202506231600"),
    @("范式进化投资组合", "513100", "纳指ETF", 1, 636.5814125697908, 1.56, 993.0670036088737, 99639.24647381395, "202506231600"),
    @("范式进化投资组合", "513290", "纳指生物科技ETF", 0.98, 884.4538210040457, 1.1, 972.8992031044504, 99639.24647381395, "202506231600"),
    @("范式进化投资组合", "518880", "黄金ETF", 0.99, 132.2001081659486, 7.46, 986.2128069179765, 99639.24647381395, "202506231600"),
    @("范式进化投资组合", "600085", "同仁堂", 0.99, 27.49471300507762, 36, 989.8096681827943, 99639.24647381395, "202506231600"),
    @("范式进化投资组合", "600900", "长江电力", 1.01, 32.81132034584936, 30.65, 1005.666968600283, 99639.24647381395, "202506231600"),
    @("范式进化投资组合", "600989", "宝丰能源", 5.04, 302.6749902285196, 16.58, 5018.351337988855, 99639.24647381395, "202506231600"),
    @("范式进化投资组合", "601899", "XD紫金矿", 9.869999999999999, 531.0482559694855, 18.51, 9829.703217995177, 99639.24647381395, "202506231600"),
    @("范式进化投资组合", "603119", "浙江荣泰", 0.95, 24.16423640557475, 39.01, 942.646862181471, 99639.24647381395, "202506231600"),
    @("范式进化投资组合", "HK06881", "中国银河", 0.96, 117.8576436007749, 8.109999999999999, 955.8254896022844, 99639.24647381395, "202506231600"),
    @("范式进化投资组合", "100000", "现金", 57.17, 56967.67061087058, 1, 56967.67061087058, 99639.24647381395, "202506231600")
)
$startRow = 166
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = "'" + $row[8]
}
